# Add a new "SUBSTITUTE_goofy" worksheet, as the last tab, that documents
# the mojibake / "goofy" characters found in the workbook (and their
# SUBSTITUTE()-ready replacements) together with a LEN() sanity check.

$wb = $excel.ActiveWorkbook

# --- create the new sheet as the last tab -----------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "SUBSTITUTE_goofy"

# --- header row ---------------------------------------------------------
$ws.Range("A1").Value = "a"
$ws.Range("B1").Value = "b"
$ws.Range("C1").Value = "c"

# --- goofy-character / replacement pairs, with a LEN() check ----------
$ws.Range("A2").Value = "â€¦"
$ws.Range("B2").Value = "…"
$ws.Range("C2").Formula = "=LEN(A2)"

$ws.Range("A3").Value = "â€“"
$ws.Range("B3").Value = "-"
$ws.Range("C3").Formula = "=LEN(A3)"

$ws.Range("A4").Value = "â€”"
$ws.Range("B4").Value = "-"
$ws.Range("C4").Formula = "=LEN(A4)"

$ws.Range("A5").Value = "â€™"
$ws.Range("B5").Value = "’"
$ws.Range("C5").Formula = "=LEN(A5)"

$ws.Range("A6").Value = "Ã©"
$ws.Range("B6").Value = "é"
$ws.Range("C6").Formula = "=LEN(A6)"

$ws.Range("A7").Value = "Ã¼"
$ws.Range("B7").Value = "ü"
$ws.Range("C7").Formula = "=LEN(A7)"

$ws.Range("A8").Value = "–"
$ws.Range("B8").Value = "-"
$ws.Range("C8").Formula = "=LEN(A8)"

$ws.Range("A9").Value = "–"
$ws.Range("B9").Value = "-"
$ws.Range("C9").Formula = "=LEN(A9)"

$ws.Range("A10").Value = "–"
$ws.Range("B10").Value = "-"
$ws.Range("C10").Formula = "=LEN(A10)"

$ws.Range("A11").Value = "—"
$ws.Range("B11").Value = "-"
$ws.Range("C11").Formula = "=LEN(A11)"

$ws.Range("A12").Value = "…"
$ws.Range("B12").Value = "…"
$ws.Range("C12").Formula = "=LEN(A12)"

# --- formatting: columns A and B down to row 199 are text-formatted ---
$ws.Range("A2:B199").NumberFormat = "@"

# --- select the whole sheet (mirrors the authored file) and make this
#     the active tab ------------------------------------------------------
$ws.Activate()
$ws.Cells.Select()
